# Hotfix: Thu Nov  7 17:45:45 RTZ 2024
# Adds a new row to the "Links" sheet and five new rows to the "Bash" sheet.

$wb = $excel.ActiveWorkbook

# ---- Links sheet: append row 17 ----
$links = $wb.Worksheets.Item("Links")
$links.Cells.Item(17, 1).Value = 52
$links.Cells.Item(17, 2).Value = "Ссылка на очень удобную программу Samsung Dex. `nПрограмма предназначена на для трансляции экрана телефона Samsung на ПК."
$links.Cells.Item(17, 3).Value = "https://www.samsung.com/ru/apps/samsung-dex/"

# ---- Bash sheet: append rows 76-80 ----
$bash = $wb.Worksheets.Item("Bash")

$bash.Cells.Item(76, 1).Value = 129
$bash.Cells.Item(76, 2).Value = "~/AppData/Local/Programs/Python/Python313/python.exe -m webbrowser http://127.0.0.1:82"
$bash.Cells.Item(76, 3).Value = "Запуск приложения в браузере"

$bash.Cells.Item(77, 1).Value = 130
$bash.Cells.Item(77, 2).Value = "Set-ExecutionPolicy RemoteSigned"
$bash.Cells.Item(77, 3).Value = "Разрешение запуска сценариев powershell"

$bash.Cells.Item(78, 1).Value = 131
$bash.Cells.Item(78, 2).Value = "where powershell"
$bash.Cells.Item(78, 3).Value = "Как найти исполняемый файл powershell в системе?`nВвести в cmd команду!"

$bash.Cells.Item(79, 1).Value = 132
$bash.Cells.Item(79, 2).Value = "python -m pip uninstall -r requirements.txt -y"
$bash.Cells.Item(79, 3).Value = "Принудительное деинсталляция пакетов из файла requirements.txt "

$bash.Cells.Item(80, 1).Value = 133
$bash.Cells.Item(80, 2).Value = "#!/c/Windows/System32/WindowsPowerShell/v1.0/powershell.exe`nStop-Process -Name `"python`""
$bash.Cells.Item(80, 3).Value = "Скрипт powershell, который убивает все python запущенные процессы"
